$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.151055335998535
$ws.Range("B1").Value = 2.59734582901001
$ws.Range("C1").Value = 3.501146078109741
$ws.Range("D1").Value = 6.229983329772949
$ws.Range("E1").Value = 1.969773292541504
